$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.550.05'
$ws.Range("E2").Value = '  -0.35%  '
$ws.Range("D3").Value = '2.509.48'
$ws.Range("E3").Value = '  -1.76%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.47'
$ws.Range("E5").Value = '  +2.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.71'
$ws.Range("E6").Value = '  -1.92%  '
$ws.Range("E7").Value = '  +2.05%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.537'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.40'
$ws.Range("E10").Value = '  +1.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0811'
$ws.Range("E11").Value = '  +0.49%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.59'
$ws.Range("E12").Value = '  +1.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.112'
$ws.Range("E13").Value = '  -3.74%  '
$ws.Range("D14").Value = '2.895.54'
$ws.Range("E14").Value = '  -1.80%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.62'
$ws.Range("E15").Value = '  +8.04%  '
$ws.Range("D16").Value = '2.532.63'
$ws.Range("E16").Value = '  -2.79%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.855'
$ws.Range("E17").Value = '  -2.82%  '
$ws.Range("D18").Value = '42.523.57'
$ws.Range("E18").Value = '  -0.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.85'
$ws.Range("E19").Value = '  -4.73%  '
$ws.Range("E20").Value = '  -1.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.51'
$ws.Range("E21").Value = '  -1.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.45'
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '252.95'
$ws.Range("E23").Value = '  -1.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.93'
$ws.Range("E24").Value = '  -0.31%  '
$ws.Range("E25").Value = '  -2.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '27.00'
$ws.Range("E26").Value = '  -3.81%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.998'
$ws.Range("E27").Value = '  -0.21%  '
$ws.Range("E28").Value = '  +10.79%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.16'
$ws.Range("E29").Value = '  +1.36%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '37.58'
$ws.Range("E30").Value = '  -3.33%  '
$ws.Range("E31").Value = '  -1.38%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '153.35'
$ws.Range("E32").Value = '  -1.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.21'
$ws.Range("E33").Value = '  +4.75%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.28'
$ws.Range("E34").Value = '  -1.13%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0785'
$ws.Range("E35").Value = '  -2.09%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.07'
$ws.Range("E36").Value = '  -4.55%  '
$ws.Range("E37").Value = '  -4.71%  '
$ws.Range("E38").Value = '  -1.48%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '24.83'
$ws.Range("E39").Value = '  -4.59%  '
$ws.Range("E40").Value = '  +0.54%  '
$ws.Range("E41").Value = '  +0.71%  '
$ws.Range("E42").Value = '  +0.40%  '
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0302'
$ws.Range("E44").Value = '  -0.56%  '
$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.01'
$ws.Range("E45").Value = '  -2.65%  '
$ws.Range("D46").Value = '2.025.76'
$ws.Range("E46").Value = '  -1.62%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '84.28'
$ws.Range("E47").Value = '  -4.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.92'
$ws.Range("E48").Value = '  -3.10%  '
$ws.Range("D49").Value = '2.753.20'
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.63'
$ws.Range("E51").Value = '  -4.45%  '
